# Append a new product row (row 23) to the ProductInfo sheet, duplicating
# the formatting of the existing data rows and reusing row 5's product data
# (GIÀY BOOTS NỮ PIERRE CARDIN CLOUDY PCWFWS153) with a new sequential ID.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductInfo")

$newRow = 23

# Copy the formatting (fill/alignment/etc.) of the row directly above so the
# new row matches the rest of the table, then overwrite with the new values.
$ws.Range("A22:E22").Copy() | Out-Null
$ws.Range("A23:E23").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item($newRow, 1).Value = "22"
$ws.Cells.Item($newRow, 2).Value = "GIÀY BOOTS NỮ PIERRE CARDIN CLOUDY PCWFWS153"
$ws.Cells.Item($newRow, 3).Value = "98.000 VND"
$ws.Cells.Item($newRow, 4).Value = "0 Đánh giá"
$ws.Cells.Item($newRow, 5).Value = "0 Đã bán"
